$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 02:47"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6635933
$ws.Range("C4").Value = 46286
$ws.Range("D4").Value = 3916264
$ws.Range("E4").Value = 2522274
$ws.Range("G4").Value = 1068
$ws.Range("H4").Value = 197395

# Row 32 - Ecuador
$ws.Range("B32").Value = 114732
$ws.Range("C32").Value = 1526
$ws.Range("E32").Value = 12654
$ws.Range("G32").Value = 87
$ws.Range("H32").Value = 10836

# Row 37 - Panama
$ws.Range("B37").Value = 100330
$ws.Range("C37").Value = 615
$ws.Range("D37").Value = 72858
$ws.Range("E37").Value = 25332
$ws.Range("G37").Value = 13
$ws.Range("H37").Value = 2140

# Row 96 - Guinea
$ws.Range("B96").Value = 9979
$ws.Range("C96").Value = 33
$ws.Range("D96").Value = 9189
$ws.Range("E96").Value = 727

# Row 102 - Gabon
$ws.Range("B102").Value = 8643
$ws.Range("C102").Value = 22
$ws.Range("D102").Value = 7706
$ws.Range("E102").Value = 884

# Row 123 - Surinam
$ws.Range("B123").Value = 4529
$ws.Range("C123").Value = 52
$ws.Range("D123").Value = 3747
$ws.Range("E123").Value = 689

# Rows 163-169: countries updated and re-sorted by total cases.
# Lesoto (was row 164) moves up to row 163 with fresh data;
# Niger (was row 163) moves down to row 164 keeping its previous data.
$ws.Range("A163").Value = "Lesoto"
$ws.Range("B163").Value = 1245
$ws.Range("C163").Value = 81
$ws.Range("D163").Value = 528
$ws.Range("E163").Value = 684
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 33

$ws.Range("A164").Value = "Niger"
$ws.Range("B164").Value = 1178
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 1100
$ws.Range("E164").Value = 9
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 69

# Polinesia Francesa (was row 169) moves up to row 167 with fresh data;
# Martinica (was row 167) moves down to row 168 keeping its previous data;
# Santo Tome y Principe (was row 168) moves down to row 169 keeping its previous data.
$ws.Range("A167").Value = "Polinesia Francesa"
$ws.Range("B167").Value = 953
$ws.Range("C167").Value = 96
$ws.Range("D167").Value = 642
$ws.Range("E167").Value = 309
$ws.Range("G167").Value = 2
$ws.Range("H167").Value = 2

$ws.Range("A168").Value = "Martinica"
$ws.Range("B168").Value = 939
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 98
$ws.Range("E168").Value = 823
$ws.Range("H168").Value = 18

$ws.Range("A169").Value = "Santo Tome y Principe"
$ws.Range("B169").Value = 906
$ws.Range("C169").Value = 5
$ws.Range("D169").Value = 866
$ws.Range("E169").Value = 25
$ws.Range("H169").Value = 15

# Row 175 - Papua Nueva Guinea
$ws.Range("B175").Value = 508
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 271

# Row 198 - San Vicente y las Granadinas
$ws.Range("D198").Value = 61
$ws.Range("E198").Value = 1
